{"js": "// Update the date heading and the 25 division problems in the table,\n// matching the author's commit (\"Update master to output generated at c986bee\").\n\n// 1) Update the date paragraph (first paragraph of the body).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.insertText(\"2024-09-16 Monday\", Word.InsertLocation.replace);\n\n// 2) Update each division-problem cell in the single table.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Only every 4th row (0, 4, 8, 12, 16) holds the 5 problems per line;\n// the rows in between are blank spacer rows left untouched.\nconst newValues = [\n  [\"23\u00f74=\", \"82\u00f78=\", \"89\u00f78=\", \"82\u00f77=\", \"27\u00f73=\"],\n  [\"71\u00f75=\", \"71\u00f78=\", \"38\u00f72=\", \"60\u00f77=\", \"91\u00f73=\"],\n  [\"25\u00f73=\", \"79\u00f74=\", \"28\u00f76=\", \"80\u00f72=\", \"92\u00f79=\"],\n  [\"16\u00f79=\", \"22\u00f78=\", \"62\u00f77=\", \"80\u00f77=\", \"60\u00f78=\"],\n  [\"94\u00f75=\", \"70\u00f77=\", \"69\u00f74=\", \"73\u00f76=\", \"56\u00f74=\"],\n];\n\nconst cellParagraphs = [];\nfor (let i = 0; i < newValues.length; i++) {\n  const rowIndex = i * 4;\n  for (let col = 0; col < newValues[i].length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    const firstPara = cell.body.paragraphs.getFirst();\n    cellParagraphs.push({ para: firstPara, text: newValues[i][col] });\n  }\n}\n\nawait context.sync();\n\nfor (const { para, text } of cellParagraphs) {\n  para.insertText(text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division problems in the table,\n# matching the author's commit (\"Update master to output generated at c986bee\").\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph (first paragraph of the document).\n$d.Paragraphs(1).Range.Text = \"2024-09-16 Monday\"\n\n# 2) Update each division-problem cell in the single table.\n# Only every 4th row (1, 5, 9, 13, 17 in 1-based indexing) holds the 5\n# problems per line; the rows in between are blank spacer rows left\n# untouched.\n$t = $d.Tables(1)\n\n$newValues = @(\n    @(\"23\u00f74=\", \"82\u00f78=\", \"89\u00f78=\", \"82\u00f77=\", \"27\u00f73=\"),\n    @(\"71\u00f75=\", \"71\u00f78=\", \"38\u00f72=\", \"60\u00f77=\", \"91\u00f73=\"),\n    @(\"25\u00f73=\", \"79\u00f74=\", \"28\u00f76=\", \"80\u00f72=\", \"92\u00f79=\"),\n    @(\"16\u00f79=\", \"22\u00f78=\", \"62\u00f77=\", \"80\u00f77=\", \"60\u00f78=\"),\n    @(\"94\u00f75=\", \"70\u00f77=\", \"69\u00f74=\", \"73\u00f76=\", \"56\u00f74=\")\n)\n\nfor ($i = 0; $i -lt $newValues.Length; $i++) {\n    $rowIndex = $i * 4 + 1\n    $rowValues = $newValues[$i]\n    for ($col = 0; $col -lt $rowValues.Length; $col++) {\n        $t.Cell($rowIndex, $col + 1).Range.Text = $rowValues[$col]\n    }\n}\n"}
